$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.955.63'
$ws.Range("E2").Value = '  -5.27%  '

$ws.Range("D3").Value = '1.821.28'
$ws.Range("E3").Value = '  -4.57%  '

$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -1.09%  '

$ws.Range("D5").Value = "'326.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.26%  '

$ws.Range("E6").Value = '  -0.79%  '

$ws.Range("D7").Value = "'0.4628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.41%  '

$ws.Range("D8").Value = "'0.3841"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.58%  '

$ws.Range("E9").Value = '  -3.71%  '

$ws.Range("D10").Value = "'0.07825"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.88%  '

$ws.Range("D11").Value = "'0.9577"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.63%  '

$ws.Range("D12").Value = "'21.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.43%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.799.28'
$ws.Range("E13").Value = '  -6.92%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = "'5.644"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.29%  '

$ws.Range("D15").Value = "'6.849"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.28%  '

$ws.Range("D16").Value = "'0.06873"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.51%  '

$ws.Range("D17").Value = "'1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.03%  '

$ws.Range("D18").Value = "'86.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.30%  '

$ws.Range("D19").Value = "'0.000009925"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.10%  '

$ws.Range("D20").Value = "'16.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.02%  '

$ws.Range("D21").Value = "'1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.53%  '

$ws.Range("D22").Value = '27.999.85'
$ws.Range("E22").Value = '  -5.23%  '

$ws.Range("D23").Value = "'5.322"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.68%  '

$ws.Range("D24").Value = "'10.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.61%  '

$ws.Range("D25").Value = "'2.127"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.39%  '

$ws.Range("D26").Value = '2.066.12'
$ws.Range("E26").Value = '  -5.16%  '

$ws.Range("D27").Value = "'151.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.21%  '

$ws.Range("D28").Value = "'19.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.05%  '

$ws.Range("D29").Value = "'5.700"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -12.50%  '

$ws.Range("D30").Value = "'1.968"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.25%  '

$ws.Range("D31").Value = "'116.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.91%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = "'0.9366"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.68%  '

$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").Value = "'0.09256"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.11%  '

$ws.Range("D34").Value = "'5.280"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.03%  '

$ws.Range("D35").Value = "'3.423"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.37%  '

$ws.Range("D36").Value = "'1.309"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.79%  '

$ws.Range("D37").Value = "'0.05943"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.44%  '

$ws.Range("D38").Value = "'0.02141"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.69%  '

$ws.Range("D39").Value = "'1.147"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.83%  '

$ws.Range("D40").Value = "'1.003"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.57%  '

$ws.Range("D41").Value = "'7.557"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.88%  '

$ws.Range("D42").Value = "'0.5583"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.25%  '

$ws.Range("D43").Value = "'9.886"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.52%  '

$ws.Range("D44").Value = "'0.1767"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.94%  '

$ws.Range("D45").Value = "'1.239"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.90%  '

$ws.Range("D46").Value = "'2.241"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.50%  '

$ws.Range("D47").Value = "'11.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.61%  '

$ws.Range("D48").Value = "'0.5248"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.12%  '

$ws.Range("D49").Value = "'0.07001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.79%  '

$ws.Range("D50").Value = "'1.820"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.73%  '

$ws.Range("D51").Value = "'112.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.08%  '
